$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the _GoBack bookmark that currently sits right after "9223"
#    (it gets re-created later, between " your appli" and "cation on ").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) ", Weather, Messaging, or Productivity " ->
#    ", Weather, Messaging, Productivity, Betting, or Courses sign up "
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    ", Weather, Messaging, or Productivity ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    ", Weather, Messaging, Productivity, Betting, or Courses sign up ", 2) | Out-Null

# ------------------------------------------------------------------
# 3) ". Below you can find a project breakdown" ->
#    ". Below you can find project breakdown"
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    ". Below you can find a project breakdown", $true, $false, $false,
    $false, $false, $true, 1, $false,
    ". Below you can find project breakdown", 2) | Out-Null

# ------------------------------------------------------------------
# 4) "Deploy your application on " ->
#       "Host" | " your appli" | <_GoBack bookmark> | "cation on "
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "Deploy your application on ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Host your application on ", 2) | Out-Null
# $rng now spans the freshly written "Host your application on " text.

# Force a run break between "Host" and " your appli" by briefly planting a
# bookmark at that offset and removing it again; the split survives even
# though the bookmark itself is gone.
$splitPoint = $d.Range($rng.Start + 4, $rng.Start + 4)
$d.Bookmarks.Add("_TMP_SPLIT", $splitPoint) | Out-Null
$d.Bookmarks("_TMP_SPLIT").Delete()

# Re-create the _GoBack bookmark between " your appli" and "cation on ".
$bookmarkPoint = $d.Range($rng.Start + 15, $rng.Start + 15)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint) | Out-Null

# ------------------------------------------------------------------
# 5) "Automated testing: s" ->
#       "Automated testing" (underlined) | ": s" (not underlined)
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(
    "Automated testing", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$rng2.Font.Underline = 1
